# Insert one new weekly price record as row 426 in the "Hortaliza, Feria
# Lagunitas de Puerto Montt - Ají" sheet, pushing the existing rows 426-445
# down to 427-446 (dimension grows from A1:R445 to A1:R446).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 426..445 down by one to make room for the new record.
$ws.Rows.Item(426).Insert()

# Populate the newly inserted row 426 with the new record's data.
$ws.Cells.Item(426, 1).Value = 4
$ws.Cells.Item(426, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(426, 3).Value = "Los Lagos"
$ws.Cells.Item(426, 4).Value = 45147
$ws.Cells.Item(426, 5).Value = 10
$ws.Cells.Item(426, 6).Value = 100112021
$ws.Cells.Item(426, 7).Value = "Ají"
$ws.Cells.Item(426, 8).Value = "Inferno"
$ws.Cells.Item(426, 9).Value = "Primera"
$ws.Cells.Item(426, 10).Value = 25
$ws.Cells.Item(426, 11).Value = 22000
$ws.Cells.Item(426, 12).Value = 22000
$ws.Cells.Item(426, 13).Value = 22000
$ws.Cells.Item(426, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(426, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(426, 16).Value = 2200
$ws.Cells.Item(426, 17).Value = 10
$ws.Cells.Item(426, 18).Value = "Hortaliza"
